$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "月累计门诊病人数" (monthly cumulative outpatient count) column (Q) is no
# longer needed; delete the entire column and let everything to its right
# (values, styles, comments, data validation, autofilter, defined names)
# shift left by one, which is exactly what Excel does for EntireColumn.Delete.
$ws.Range("Q1").EntireColumn.Delete()

# Restore the selection recorded for the sheet after the edit.
$ws.Range("O7").Select()
